$wb = $excel.ActiveWorkbook

# Rename the first sheet from "Meal Calender" to "Meal Calendar" (typo fix)
$calSheet = $wb.Worksheets.Item(1)
$calSheet.Name = "Meal Calendar"

# Add new "Recipes" worksheet (gets inserted before the active sheet)
$recipes = $wb.Worksheets.Add()
$recipes.Name = "Recipes"

# Move the new sheet to the end (after "Meal Calendar")
$recipes.Move(1)

# Re-resolve by name: the sheet reference can go stale (rebind by index)
# across a Move(), so fetch a fresh handle before writing to it.
$recipes = $wb.Worksheets.Item("Recipes")

# Populate the Recipes sheet with header + placeholder meal names
$recipes.Range("A1").Value = "Meal Name"
$recipes.Range("A2").Value = "Salad"
$recipes.Range("A3").Value = "Soup"
$recipes.Range("A4").Value = "Egg"
$recipes.Range("A5").Value = "Beef"
$recipes.Range("A6").Value = "Chicken"
$recipes.Range("A7").Value = "PlaceHolder"

# Make Recipes the active sheet / selection
$recipes.Select()
$recipes.Range("A7").Select()
